$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-10-04"

# Update the October row label text
$ws.Range("A11").Value = "October (through 10-04)"

# Update October row (row 11) data values
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 22
$ws.Range("H11").Value = 31
$ws.Range("I11").Value = 13

# Update Total row (row 12) data values
$ws.Range("B12").Value = 230
$ws.Range("C12").Value = 435
$ws.Range("D12").Value = 637
$ws.Range("E12").Value = 558
$ws.Range("F12").Value = 425
$ws.Range("G12").Value = 923
$ws.Range("H12").Value = 1278
$ws.Range("I12").Value = 1295
